# Removed HOL slides from slide decks
#
# The deck has a "Hands-On Lab" slide (title "Hands-On Lab", body
# "SLURM Linux Cluster HOL.html") that needs to be removed. It is the
# 8th slide in the deck (1-based, matching PowerPoint's Slides collection
# / the presentation's sldIdLst order).

$p = $ppt.ActivePresentation

$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "Hands-On Lab") {
                $target = $slide
            }
        }
    }
}

if ($target -ne $null) {
    $target.Delete()
}
